$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.508.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.566.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.40%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.564.76'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.35%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.13%  '
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.172.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.55%  '
$ws.Range('E14').Value = '  +3.91%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.569.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.417.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.14%  '
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '396.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').Value = '  +4.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.710.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +10.54%  '
$ws.Range('E28').Value = '  +8.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.583.89'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '24.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.03%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  +3.51%  '
$ws.Range('E37').Value = '  +2.47%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '168.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0802'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.835'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.62'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +15.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.44%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.93%  '
$ws.Range('E48').Value = '  +7.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.452.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.82'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.15'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.02%  '
